# Add "2022-Q1" quarterly data:
#  - the existing "总计" sheet (4th sheet, holding the running totals
#    table) becomes the new "2022-Q1" fund-holding sheet
#  - a brand-new totals sheet is appended at the end, seeded from the old
#    totals table plus a new leading row for "2022-Q1", then renamed to
#    "总计"

$wb = $excel.ActiveWorkbook

$oldTotal = $wb.Worksheets.Item(4)      # currently the totals sheet
$template  = $wb.Worksheets.Item(3)     # "2021-Q4", used as a layout template

# --- 1. Build the new totals sheet first, from the current totals table ---
$newTotal = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$oldTotal.Range("A1:D4").Copy($newTotal.Range("A1:D4"))

# Insert a fresh row for the 2022-Q1 entry, pushing the rest down
$newTotal.Rows.Item(2).Insert()

# Re-number the shifted index column (A) for the old rows
$newTotal.Cells.Item(3,1).Value = 1
$newTotal.Cells.Item(4,1).Value = 2
$newTotal.Cells.Item(5,1).Value = 3

# Give the new row the same look as the rows around it
$newTotal.Range("A3:D3").Copy($newTotal.Range("A2:D2"))

$newTotal.Cells.Item(2,1).Value = 0
$newTotal.Cells.Item(2,2).Value = "2022-Q1"
$newTotal.Cells.Item(2,3).Value = 1
$newTotal.Cells.Item(2,4).Value = 0.03

# --- 2. Turn the old totals sheet into the new "2022-Q1" holdings sheet ---
# (renamed before the new sheet takes the old name, to avoid a name clash)
$oldTotal.Name = "2022-Q1"

$oldTotal.Cells.Clear()
$template.Range("A1:H2").Copy($oldTotal.Range("A1:H2"))

$oldTotal.Range("B2").NumberFormat = "@"
$oldTotal.Range("D2:G2").NumberFormat = "@"

$oldTotal.Cells.Item(2,2).Value = "320017"
$oldTotal.Cells.Item(2,3).Value = "诺安全球收益不动产(QDII)"
$oldTotal.Cells.Item(2,4).Value = "0.29"
$oldTotal.Cells.Item(2,5).Value = "93.32"
$oldTotal.Cells.Item(2,6).Value = "8.93"
$oldTotal.Cells.Item(2,7).Value = "0.0259"
$oldTotal.Cells.Item(2,8).Value = 1

# --- 3. Finally give the new sheet the totals sheet's name ---
$newTotal.Name = "总计"
